# Scheduled-runner price/profit refresh across the Sophia_Profits sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose market data changed, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H21").Value = 23510.5
$ws.Range("I21").Value = 12000
$ws.Range("J21").Value = 35021
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 35021
$ws.Range("M21").Value = -11532
$ws.Range("N21").Value = -35957
$ws.Range("H23").Value = 23510.5
$ws.Range("I23").Value = 12000
$ws.Range("J23").Value = 35021
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 35021
$ws.Range("M23").Value = -11766
$ws.Range("N23").Value = -35489
$ws.Range("H34").Value = 8763.5
$ws.Range("I34").Value = 8334.666999999999
$ws.Range("K34").Value = 8334.666999999999
$ws.Range("M34").Value = -8131.666999999999
$ws.Range("H36").Value = 8763.5
$ws.Range("I36").Value = 8334.666999999999
$ws.Range("K36").Value = 8334.666999999999
$ws.Range("M36").Value = -7619.666999999999
$ws.Range("H53").Value = 292.75
$ws.Range("I53").Value = 288.75
$ws.Range("J53").Value = 294.75
$ws.Range("K53").Value = 288.75
$ws.Range("L53").Value = 294.75
$ws.Range("M53").Value = 348.25
$ws.Range("N53").Value = -1568.75
$ws.Range("H98").Value = 2615.6667
$ws.Range("I98").Value = 1010
$ws.Range("K98").Value = 1010
$ws.Range("M98").Value = 488
$ws.Range("H118").Value = 2525.7632
$ws.Range("J118").Value = 2999.9033
$ws.Range("L118").Value = 8999.7099
$ws.Range("N118").Value = -12313.7099
$ws.Range("H122").Value = 2615.6667
$ws.Range("I122").Value = 1010
$ws.Range("K122").Value = 3030
$ws.Range("M122").Value = -580
$ws.Range("H135").Value = 921.56525
$ws.Range("I135").Value = 757.5714
$ws.Range("J135").Value = 2643.5
$ws.Range("K135").Value = 6818.1426
$ws.Range("L135").Value = 23791.5
$ws.Range("M135").Value = -4283.1426
$ws.Range("N135").Value = -28861.5
$ws.Range("H137").Value = 2271.1667
$ws.Range("I137").Value = 2136.5715
$ws.Range("K137").Value = 6409.7145
$ws.Range("M137").Value = -3859.7145
$ws.Range("H138").Value = 5548.0356
$ws.Range("J138").Value = 6998.1
$ws.Range("L138").Value = 20994.3
$ws.Range("N138").Value = -31274.3
$ws.Range("H141").Value = 1665
$ws.Range("I141").Value = 1997.5
$ws.Range("K141").Value = 5992.5
$ws.Range("M141").Value = -812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 17559.445
$ws.Range("I74").Value = 18802.8
$ws.Range("K74").Value = 18802.8
$ws.Range("M74").Value = -17928.8
$ws.Range("H77").Value = 17559.445
$ws.Range("I77").Value = 18802.8
$ws.Range("K77").Value = 94014
$ws.Range("M77").Value = -89646
$ws.Range("H102").Value = 2948.8572
$ws.Range("J102").Value = 3399.9
$ws.Range("L102").Value = 3399.9
$ws.Range("N102").Value = -6643.9
$ws.Range("H122").Value = 4115.143
$ws.Range("I122").Value = 4555.8
$ws.Range("K122").Value = 13667.4
$ws.Range("M122").Value = -11217.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9874.375
$ws.Range("I20").Value = 15349.5
$ws.Range("K20").Value = 15349.5
$ws.Range("M20").Value = -15102.5
$ws.Range("H81").Value = 32500
$ws.Range("J81").Value = 32500
$ws.Range("L81").Value = 32500
$ws.Range("N81").Value = -34622
$ws.Range("H84").Value = 32500
$ws.Range("J84").Value = 32500
$ws.Range("L84").Value = 97500
$ws.Range("N84").Value = -108108
$ws.Range("H86").Value = 7700
$ws.Range("I86").Value = 2066.6667
$ws.Range("J86").Value = 13333.333
$ws.Range("K86").Value = 2066.6667
$ws.Range("L86").Value = 13333.333
$ws.Range("M86").Value = -943.6667000000002
$ws.Range("N86").Value = -15579.333
$ws.Range("H89").Value = 7700
$ws.Range("I89").Value = 2066.6667
$ws.Range("J89").Value = 13333.333
$ws.Range("K89").Value = 10333.3335
$ws.Range("L89").Value = 66666.66500000001
$ws.Range("M89").Value = -4717.333500000001
$ws.Range("N89").Value = -77898.66500000001
$ws.Range("H99").Value = 34403.332
$ws.Range("I99").Value = 34403.332
$ws.Range("K99").Value = 34403.332
$ws.Range("M99").Value = -32905.332
$ws.Range("H135").Value = 88152.8
$ws.Range("J135").Value = 88152.8
$ws.Range("L135").Value = 88152.8
$ws.Range("N135").Value = -98292.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3210.6667
$ws.Range("I58").Value = 3191
$ws.Range("J58").Value = 3250
$ws.Range("K58").Value = 3191
$ws.Range("L58").Value = 3250
$ws.Range("M58").Value = -2988
$ws.Range("N58").Value = -3656
$ws.Range("H132").Value = 2913.6365
$ws.Range("I132").Value = 1996.1666
$ws.Range("K132").Value = 5988.4998
$ws.Range("M132").Value = -3458.4998
$ws.Range("H136").Value = 3210.6667
$ws.Range("I136").Value = 3191
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 9573
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -7023
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 138.1875
$ws.Range("I33").Value = 176.41667
$ws.Range("K33").Value = 1058.50002
$ws.Range("M33").Value = -775.5000199999999
$ws.Range("H122").Value = 1896.5
$ws.Range("J122").Value = 2120.625
$ws.Range("L122").Value = 19085.625
$ws.Range("N122").Value = -23985.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3278
$ws.Range("I122").Value = 3283.7144
$ws.Range("K122").Value = 9851.143199999999
$ws.Range("M122").Value = -7401.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5001.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5001.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5001.5
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -5341.5
$ws.Range("H22").Value = 4858.4546
$ws.Range("I22").Value = 2730
$ws.Range("J22").Value = 8583.25
$ws.Range("K22").Value = 2730
$ws.Range("L22").Value = 8583.25
$ws.Range("M22").Value = -2435
$ws.Range("N22").Value = -9173.25
$ws.Range("H27").Value = 4858.4546
$ws.Range("I27").Value = 2730
$ws.Range("J27").Value = 8583.25
$ws.Range("K27").Value = 2730
$ws.Range("L27").Value = 8583.25
$ws.Range("M27").Value = -2623
$ws.Range("N27").Value = -8797.25
$ws.Range("H82").Value = 2537.2
$ws.Range("J82").Value = 1882.5
$ws.Range("L82").Value = 1882.5
$ws.Range("N82").Value = -2604.5
$ws.Range("H85").Value = 2537.2
$ws.Range("J85").Value = 1882.5
$ws.Range("L85").Value = 1882.5
$ws.Range("N85").Value = -4378.5
$ws.Range("H100").Value = 2266.3333
$ws.Range("I100").Value = 1900
$ws.Range("K100").Value = 1900
$ws.Range("M100").Value = -1359
$ws.Range("H122").Value = 9249.25
$ws.Range("H136").Value = 3599.2856
$ws.Range("I136").Value = 2459
$ws.Range("J136").Value = 6450
$ws.Range("K136").Value = 7377
$ws.Range("L136").Value = 19350
$ws.Range("M136").Value = -4827
$ws.Range("N136").Value = -24450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 3000000
$ws.Range("I7").Value = 3000000
$ws.Range("K7").Value = 3000000
$ws.Range("M7").Value = -2999887
$ws.Range("H100").Value = 1101.7273
$ws.Range("I100").Value = 1142.8572
$ws.Range("J100").Value = 1029.75
$ws.Range("K100").Value = 2285.7144
$ws.Range("L100").Value = 2059.5
$ws.Range("M100").Value = -1744.7144
$ws.Range("N100").Value = -3141.5
$ws.Range("H132").Value = 6998.7144
$ws.Range("I132").Value = 3248
$ws.Range("K132").Value = 9744
$ws.Range("M132").Value = -7214
$ws.Range("H136").Value = 1074.75
$ws.Range("I136").Value = 1074.75
$ws.Range("K136").Value = 3224.25
$ws.Range("M136").Value = -674.25
